{"js": "// Update the worksheet date and the 25 two-digit multiplication problems.\nconst body = context.document.body;\n\n// 1. Update the date/day heading paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-01-06 Saturday\", \"Replace\");\n\n// 2. Update the multiplication problems inside the table, row by row.\n// The table has 20 rows total: 5 \"data\" rows (each holding 5 problems)\n// separated by 3 blank rows. Replacements are applied in document order\n// (row-major), matching the order the problems were regenerated in.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst dataRowIndexes = [0, 4, 9, 14, 19];\nconst newValues = [\n  [\"72\u00d783=\", \"99\u00d752=\", \"35\u00d717=\", \"16\u00d772=\", \"20\u00d721=\"],\n  [\"59\u00d782=\", \"30\u00d735=\", \"89\u00d729=\", \"28\u00d717=\", \"87\u00d746=\"],\n  [\"96\u00d763=\", \"96\u00d746=\", \"67\u00d759=\", \"52\u00d778=\", \"43\u00d717=\"],\n  [\"23\u00d776=\", \"26\u00d712=\", \"53\u00d761=\", \"99\u00d778=\", \"37\u00d795=\"],\n  [\"13\u00d745=\", \"24\u00d729=\", \"39\u00d749=\", \"96\u00d783=\", \"44\u00d723=\"],\n];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const row = dataRowIndexes[i];\n  for (let col = 0; col < 5; col++) {\n    table.getCell(row, col).value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 two-digit multiplication problems.\n$d = $word.ActiveDocument\n\n# 1. Update the date/day heading (first paragraph in the document).\n$d.Paragraphs(1).Range.Text = \"2024-01-06 Saturday\"\n\n# 2. Update the multiplication problems inside the table, row by row.\n# The table has 20 rows total: 5 \"data\" rows (each holding 5 problems)\n# separated by 3 blank rows. Replacements are applied in document order\n# (row-major), matching the order the problems were regenerated in.\n$t = $d.Tables.Item(1)\n$dataRows = @(1, 5, 10, 15, 20)\n$newValues = @(\n  @(\"72\u00d783=\", \"99\u00d752=\", \"35\u00d717=\", \"16\u00d772=\", \"20\u00d721=\"),\n  @(\"59\u00d782=\", \"30\u00d735=\", \"89\u00d729=\", \"28\u00d717=\", \"87\u00d746=\"),\n  @(\"96\u00d763=\", \"96\u00d746=\", \"67\u00d759=\", \"52\u00d778=\", \"43\u00d717=\"),\n  @(\"23\u00d776=\", \"26\u00d712=\", \"53\u00d761=\", \"99\u00d778=\", \"37\u00d795=\"),\n  @(\"13\u00d745=\", \"24\u00d729=\", \"39\u00d749=\", \"96\u00d783=\", \"44\u00d723=\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Count; $i++) {\n  $row = $dataRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $t.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n  }\n}\n"}
